$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.262.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.52%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.906.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.17%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'307.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.60%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.03%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5235"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3784"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.52%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07263"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.28%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.70%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.9015"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.28%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08225"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +9.79%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.908.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.11%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'95.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.20%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.286"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.29%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.05%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E18").Value = "'  +2.35%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.9995"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.08%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'27.296.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.43%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.068"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.40%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.155.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.35%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").Value = "'2.305"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +10.11%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'146.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.77%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -1.71%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'18.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.93%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'114.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.12%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.983"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +6.17%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.816"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.75%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.09212"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.19%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.8066"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +7.62%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.05074"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.97%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.243"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +8.06%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.921"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.81%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.331"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.75%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.576"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.93%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.5730"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.60%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.14%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.15%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'9.080"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.88%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'FraxShare"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'6.630"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.93%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Quant"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'119.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.66%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.1519"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.42%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4844"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.32%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'10.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.37%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.00%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.617"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.26%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'37.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.41%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'63.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.27%  "
$ws.Range("E51").Style = "Normal"
